$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("V1").Value = 0.8196898509754964
$ws.Range("A2").Value = 0.89786389589448246
$ws.Range("BO2").Value = 0.99434682847654321
$ws.Range("B3").Value = 0.91799713509561598
$ws.Range("AE3").Value = 0.9615654729102362
$ws.Range("AX3").Value = 0.97646705389082944
$ws.Range("BC3").Value = 0.9085873207134012
$ws.Range("C4").Value = 0.81644666083881201
$ws.Range("BJ4").Value = 0.78184717569680529
$ws.Range("J5").Value = 0.88436780936657122
$ws.Range("U5").Value = 0.95659629227246978
$ws.Range("BH5").Value = 0.58237620540852997
$ws.Range("M6").Value = 0.64776060777357736
$ws.Range("Q7").Value = 0.88647154919835902
$ws.Range("I8").Value = 0.85168623949150302
$ws.Range("D9").Value = 0.73638063343970839
$ws.Range("BL9").Value = 0.83025528862996234
$ws.Range("H10").Value = 0.55335974691888523
$ws.Range("AG10").Value = 0.85229449432408111
$ws.Range("Z11").Value = 0.706435049448629
$ws.Range("AQ11").Value = 0.73753827330045929
$ws.Range("BG11").Value = 0.60339811295594648
$ws.Range("G12").Value = 0.82531395020077225
$ws.Range("AZ12").Value = 0.83883550930200024
$ws.Range("L13").Value = 0.7778930298093778
$ws.Range("S13").Value = 0.86240126345043611
$ws.Range("T13").Value = 0.68349295551521039
$ws.Range("AV14").Value = 0.93029171453895743
$ws.Range("P15").Value = 0.8794884662952811
$ws.Range("Q15").Value = 0.92877604213375409
$ws.Range("AM17").Value = 0.75239941064087301
$ws.Range("AV17").Value = 0.96496821720182313
$ws.Range("BA18").Value = 0.77150267296377195
$ws.Range("BK18").Value = 0.84604578760754956
$ws.Range("D19").Value = 0.72733351288742842
$ws.Range("X20").Value = 0.76313759598868713
$ws.Range("BJ20").Value = 0.68445972444055414
$ws.Range("AB21").Value = 0.80787074352202648
$ws.Range("AO21").Value = 0.95453049583342464
$ws.Range("B22").Value = 0.71628889408992413
$ws.Range("U22").Value = 0.85462142332768187
$ws.Range("BC22").Value = 0.90915879405867162
$ws.Range("BM22").Value = 0.60718659411258258
$ws.Range("AH24").Value = 0.77976657115063863
$ws.Range("BD24").Value = 0.66242802628618858
$ws.Range("A25").Value = 0.97229334201554951
$ws.Range("E26").Value = 0.95723082137572724
$ws.Range("G26").Value = 0.79277668789574529
$ws.Range("BF26").Value = 0.95167924319206254
$ws.Range("V27").Value = 0.99750206289537724
$ws.Range("AS27").Value = 0.74389670322354706
$ws.Range("BD27").Value = 0.92775372743357098
$ws.Range("AG28").Value = 0.87056683442339722
$ws.Range("AY29").Value = 0.8987808906201098
$ws.Range("BD29").Value = 0.98998171252767886
$ws.Range("BF29").Value = 0.94245121925840603
$ws.Range("AB30").Value = 0.75272695439639636
$ws.Range("AE30").Value = 0.59334742979497301
$ws.Range("AT30").Value = 0.94494140393192039
$ws.Range("AU30").Value = 0.95880657416028492
$ws.Range("AB31").Value = 0.73636740690543423
$ws.Range("K32").Value = 0.79772565525357897
$ws.Range("Q32").Value = 0.92222907882404859
$ws.Range("AH32").Value = 0.82720258695762561
$ws.Range("B33").Value = 0.79595010826066948
$ws.Range("N33").Value = 0.68670457777987726
$ws.Range("P33").Value = 0.56122521446195739
$ws.Range("AA33").Value = 0.87705532647328377
$ws.Range("AC33").Value = 0.9670199912299402
$ws.Range("AJ35").Value = 0.54269208145541925
$ws.Range("AH36").Value = 0.92354848379769572
$ws.Range("AV37").Value = 0.90552698301515688
$ws.Range("BK37").Value = 0.92373206660786977
$ws.Range("AA38").Value = 0.6409977241055147
$ws.Range("AM38").Value = 0.77704789316527689
$ws.Range("AJ39").Value = 0.98182494251537666
$ws.Range("AK39").Value = 0.81399818070346341
$ws.Range("BC39").Value = 0.95858549877784771
$ws.Range("AM40").Value = 0.85243716855873775
$ws.Range("AS40").Value = 0.74568284631987014
$ws.Range("AW40").Value = 0.75004511173327848
$ws.Range("BB40").Value = 0.89108746484566426
$ws.Range("AN42").Value = 0.68287798665243926
$ws.Range("BB42").Value = 0.69513665165438554
$ws.Range("AS43").Value = 0.65451310185986877
$ws.Range("AM44").Value = 0.79009375156416084
$ws.Range("AT44").Value = 0.86734101045214917
$ws.Range("AX44").Value = 0.96492117605348815
$ws.Range("BA44").Value = 0.83991989068173212
$ws.Range("BF44").Value = 0.88617892536962772
$ws.Range("W45").Value = 0.9537980267474977
$ws.Range("X45").Value = 0.88024156490408623
$ws.Range("AB46").Value = 0.60845154581077654
$ws.Range("BP46").Value = 0.91846569110523379
$ws.Range("BG47").Value = 0.65941778683574381
$ws.Range("J48").Value = 0.62945298066798405
$ws.Range("F49").Value = 0.60275293215027137
$ws.Range("N49").Value = 0.64870187369516796
$ws.Range("Z49").Value = 0.9497985510331588
$ws.Range("AP49").Value = 0.70592969838118647
$ws.Range("AQ51").Value = 0.67054698961523362
$ws.Range("BL51").Value = 0.94849404062659803
$ws.Range("BP51").Value = 0.86602748561211773
$ws.Range("M52").Value = 0.82455974166759449
$ws.Range("X52").Value = 0.94370645727845459
$ws.Range("Y53").Value = 0.78865989698912786
$ws.Range("AI53").Value = 0.65084094823479688
$ws.Range("AJ53").Value = 0.62420351415887454
$ws.Range("AO53").Value = 0.98908867627690422
$ws.Range("R54").Value = 0.6034740791317752
$ws.Range("W54").Value = 0.95333671099144746
$ws.Range("Z54").Value = 0.99636182221272696
$ws.Range("BA54").Value = 0.85961159228116979
$ws.Range("AZ55").Value = 0.91359555109714208
$ws.Range("BA55").Value = 0.98530898612481244
$ws.Range("M56").Value = 0.84700938276899063
$ws.Range("T57").Value = 0.70546905966859197
$ws.Range("Z57").Value = 0.73376350624965359
$ws.Range("BD57").Value = 0.85135151600256087
$ws.Range("AQ58").Value = 0.97096864115692894
$ws.Range("C60").Value = 0.87818374588137416
$ws.Range("D60").Value = 0.97878089041516947
$ws.Range("S60").Value = 0.63047938228128964
$ws.Range("AU60").Value = 0.75770611284667799
$ws.Range("AM61").Value = 0.86583727083070561
$ws.Range("BC61").Value = 0.75708677417130033
$ws.Range("AW62").Value = 0.74029975325844855
$ws.Range("BH63").Value = 0.90383328365371374
$ws.Range("T64").Value = 0.74087158682203114
$ws.Range("BN64").Value = 0.93453656625460824
$ws.Range("AC65").Value = 0.95678028264902593
$ws.Range("P66").Value = 0.98548720081153118
$ws.Range("AJ66").Value = 0.91273118367834027
$ws.Range("AX66").Value = 0.92908737496532212
$ws.Range("R67").Value = 0.82922681066217596
$ws.Range("G68").Value = 0.89846351140704606

# Column width adjustments (best-effort; COM ColumnWidth rounds to 1/6-character granularity
# so the exact stored XML widths 12.7109375 / 11.7109375 cannot be hit precisely).
$ws.Columns.Item(23).ColumnWidth = 11.8
$ws.Columns.Item(25).ColumnWidth = 10.8
